$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 3.041
$ws.Range("C2").Value = 21.59
$ws.Range("D2").Value = 90.587
$ws.Range("E2").Value = 241.919
$ws.Range("F2").Value = 459.373
$ws.Range("G2").Value = 804.925
$ws.Range("H2").Value = 2106.572
$ws.Range("I2").Value = 8517.050999999999
$ws.Range("J2").Value = 12356.957
$ws.Range("K2").Value = 18223.146
$ws.Range("L2").Value = 28473.934
$ws.Range("M2").Value = 33703.88
$ws.Range("N2").Value = 73387.905

# Row 3
$ws.Range("B3").Value = 1.222
$ws.Range("C3").Value = 2.74
$ws.Range("D3").Value = 3.974
$ws.Range("E3").Value = 7.828
$ws.Range("F3").Value = 10.384
$ws.Range("G3").Value = 12.641
$ws.Range("H3").Value = 37.459
$ws.Range("I3").Value = 73.06399999999999
$ws.Range("J3").Value = 84.658
$ws.Range("K3").Value = 127.222
$ws.Range("L3").Value = 140.433
$ws.Range("M3").Value = 169.978
$ws.Range("N3").Value = 241.43

# Row 4
$ws.Range("B4").Value = 1.543
$ws.Range("C4").Value = 4.479
$ws.Range("D4").Value = 10.483
$ws.Range("E4").Value = 16.908
$ws.Range("F4").Value = 23.068
$ws.Range("G4").Value = 31.06
$ws.Range("H4").Value = 69.398
$ws.Range("I4").Value = 206.932
$ws.Range("J4").Value = 202.412
$ws.Range("K4").Value = 290.046
$ws.Range("L4").Value = 316.534
$ws.Range("M4").Value = 389.566
$ws.Range("N4").Value = 584.6

# Row 5
$ws.Range("B5").Value = 0.209
$ws.Range("C5").Value = 0.412
$ws.Range("D5").Value = 0.518
$ws.Range("E5").Value = 0.539
$ws.Range("F5").Value = 0.873
$ws.Range("G5").Value = 6.104
$ws.Range("H5").Value = 3.231
$ws.Range("I5").Value = 8.733000000000001
$ws.Range("J5").Value = 12.706
$ws.Range("K5").Value = 14.062
$ws.Range("L5").Value = 19.19
$ws.Range("M5").Value = 21.531
$ws.Range("N5").Value = 27.872

# Row 6
$ws.Range("B6").Value = 0.4
$ws.Range("C6").Value = 2.927
$ws.Range("D6").Value = 6.755
$ws.Range("E6").Value = 10.262
$ws.Range("F6").Value = 12.015
$ws.Range("G6").Value = 11.775
$ws.Range("H6").Value = 22.814
$ws.Range("I6").Value = 51.243
$ws.Range("J6").Value = 64.962
$ws.Range("K6").Value = 94.54900000000001
$ws.Range("L6").Value = 120.797
$ws.Range("M6").Value = 109.447
$ws.Range("N6").Value = 233.377

